$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates ---

# D-column values that are unambiguous text (contain a second "." or other
# non-numeric characters) can be assigned directly without Excel coercing
# them into a floating-point number.
$dTextValues = @{
    2 = "64.950.77"
    3 = "3.388.13"
    8 = "3.376.01"
    15 = "3.930.04"
    18 = "3.392.05"
    19 = "64.897.50"
    40 = "0.0₃0749"
    42 = "3.091.56"
}
foreach ($row in $dTextValues.Keys) {
    $ws.Range("D$row").Value = $dTextValues[$row]
}

# D-column values that look like plain decimal numbers (single ".") would be
# auto-coerced to a numeric type by Excel, which also introduces binary
# floating-point rounding noise. Force the cell to Text format first, assign
# the literal string, then restore the default "Normal" style so no stray
# formatting is left behind.
$dNumericTextValues = @{
    5 = "558.73"
    6 = "173.23"
    11 = "0.630"
    12 = "54.27"
    14 = "9.10"
    16 = "18.26"
    20 = "11.81"
    21 = "0.992"
    22 = "470.90"
    25 = "86.95"
    26 = "13.57"
    28 = "10.75"
    29 = "8.77"
    30 = "30.92"
    31 = "6.68"
    32 = "11.48"
    33 = "573.03"
    34 = "61.48"
    39 = "35.75"
    43 = "1.00"
    47 = "2.46"
    50 = "139.46"
    51 = "8.31"
}
foreach ($row in $dNumericTextValues.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $dNumericTextValues[$row]
    $cell.Style = "Normal"
}

# --- Column E (Volume(1h)) updates ---
# These values always carry surrounding spaces and a "%" sign, so Excel never
# mistakes them for numbers and a direct assignment is safe.
$eValues = @{
    2 = "  +2.23%  "
    3 = "  +2.10%  "
    4 = "  +0.03%  "
    5 = "  +1.99%  "
    6 = "  +0.73%  "
    7 = "  +1.59%  "
    8 = "  +2.00%  "
    9 = "  -0.01%  "
    10 = "  +11.65%  "
    11 = "  +3.14%  "
    12 = "  +2.16%  "
    13 = "  +5.43%  "
    14 = "  +2.93%  "
    15 = "  +2.11%  "
    16 = "  +1.16%  "
    17 = "  +2.12%  "
    18 = "  +1.91%  "
    19 = "  +2.07%  "
    21 = "  +2.06%  "
    22 = "  +14.40%  "
    23 = "  +13.54%  "
    24 = "  +2.55%  "
    26 = "  -0.67%  "
    27 = "  +6.69%  "
    28 = "  +2.32%  "
    29 = "  +2.11%  "
    30 = "  +6.67%  "
    31 = "  +5.10%  "
    32 = "  +1.48%  "
    33 = "  -0.11%  "
    34 = "  +6.80%  "
    35 = "  +2.18%  "
    36 = "  +0.01%  "
    37 = "  +5.23%  "
    38 = "  -4.96%  "
    39 = "  +2.12%  "
    40 = "  +1.95%  "
    41 = "  +1.48%  "
    42 = "  -0.69%  "
    43 = "  +0.13%  "
    44 = "  +2.88%  "
    45 = "  +4.10%  "
    46 = "  +5.40%  "
    47 = "  +2.40%  "
    48 = "  -2.32%  "
    49 = "  +0.23%  "
    50 = "  +5.27%  "
    51 = "  +4.07%  "
}
foreach ($row in $eValues.Keys) {
    $ws.Range("E$row").Value = $eValues[$row]
}

"Updated $($dTextValues.Count + $dNumericTextValues.Count) price cells and $($eValues.Count) volume cells"
